$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.342636585235596
$ws.Range("B1").Value = 4.472248077392578
$ws.Range("C1").Value = 8.236617088317871
$ws.Range("D1").Value = 8.498508453369141
$ws.Range("E1").Value = 5.36424446105957
